$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set E2:E9 to 1 (new column of values added alongside existing data rows)
$ws.Range("E2:E9").Value = 1

# Update the active selection on the sheet (was F17, now I10)
$ws.Range("I10").Select()
